$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Update column A values (rows 2-12) from "24" to "2"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "2"
}

# Delete column M (the "Event" header and its "nan" values)
$ws.Columns.Item(13).Delete()
